# feat: add 2022-Q3 data
#
# 1) "总计" (summary) sheet gets a new leading row for 2022-Q3, pushing the
#    existing quarters down by one.
# 2) A brand-new worksheet "2022-Q3" is inserted right after "总计" (so it
#    becomes the 2nd tab), holding the per-fund holdings detail for the
#    quarter. The sheets that follow keep their original order/content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 6 is brand new (previously the sheet only went to row 5) - clone the
# bold/bordered "index column" style from row 5's A cell so A6 matches the
# look of A2:A5 instead of picking up the default style.
$summary.Cells.Item(5, 1).Copy($summary.Cells.Item(6, 1))

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 36
$summary.Cells.Item(2, 4).Value = 4.56

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q2"
$summary.Cells.Item(3, 3).Value = 20
$summary.Cells.Item(3, 4).Value = 3.6

$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2022-Q1"
$summary.Cells.Item(4, 3).Value = 5
$summary.Cells.Item(4, 4).Value = 0.6

$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(5, 2).Value = "2021-Q4"
$summary.Cells.Item(5, 3).Value = 6
$summary.Cells.Item(5, 4).Value = 1.28

$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(6, 2).Value = "2021-Q3"
$summary.Cells.Item(6, 3).Value = 10
$summary.Cells.Item(6, 4).Value = 2.93

# ---------------------------------------------------------------------
# Part 2: new "2022-Q3" worksheet, inserted right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Match the page-margin convention the rest of the workbook uses (values
# are points: 0.75in/0.75in/1in/1in/0.5in/0.5in).
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# This sheet uses the exact same header/index-column look as every other
# quarterly detail sheet in the workbook - clone that formatting from the
# "2022-Q2" sheet (bold header row, bordered/bold column A) instead of
# re-deriving it by hand, so the style table stays aligned with theirs.
$template = $wb.Worksheets.Item("2022-Q2")
$template.Range("B1:H1").Copy($q3.Range("B1:H1"))
for ($i = 2; $i -le 37; $i++) {
    $template.Cells.Item(2, 1).Copy($q3.Cells.Item($i, 1))
}

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking *text* in this workbook (fund
# codes with leading zeros, and figures kept as formatted strings) - force
# text storage before writing so Excel doesn't silently coerce them to
# numbers (which would e.g. drop the leading zero on fund code "050004").
$q3.Range("B2:B37").NumberFormat = "@"
$q3.Range("D2:G37").NumberFormat = "@"

$data = @(
  @(0, "481010", "工银中小盘混合", "15.67", "87.16", "4.11", "0.6440", 2),
  @(1, "050004", "博时精选混合A", "19.78", "66.67", "2.93", "0.5796", 7),
  @(2, "012985", "平安优势回报1年持有期混合型证券投资基金A", "13.50", "93.07", "3.56", "0.4806", 7),
  @(3, "481015", "工银主题策略混合A", "8.53", "90.81", "4.72", "0.4026", 5),
  @(4, "210003", "金鹰行业优势混合", "5.59", "82.52", "5.98", "0.3343", 1),
  @(5, "013417", "博时核心资产精选混合A", "7.47", "70.34", "3.31", "0.2473", 6),
  @(6, "012917", "平安优势领航1年持有期混合A", "7.23", "93.10", "3.39", "0.2451", 8),
  @(7, "011884", "工银景气优选混合A", "6.34", "86.74", "2.98", "0.1889", 9),
  @(8, "002450", "平安睿享文娱灵活配置混合A", "3.64", "88.33", "4.12", "0.1500", 5),
  @(9, "011260", "金鹰新能源混合A", "3.01", "88.11", "4.38", "0.1318", 10),
  @(10, "010126", "平安价值成长混合A", "3.36", "93.06", "3.63", "0.1220", 6),
  @(11, "006101", "平安优势产业灵活配置混合C", "3.03", "93.46", "3.56", "0.1079", 8),
  @(12, "011261", "金鹰新能源混合C", "2.44", "88.11", "4.38", "0.1069", 10),
  @(13, "011828", "平安睿享成长混合A", "2.57", "93.00", "4.07", "0.1046", 6),
  @(14, "501063", "汇添富悦享定期开放混合", "2.19", "66.28", "3.34", "0.0731", 7),
  @(15, "002451", "平安睿享文娱灵活配置混合C", "1.72", "88.33", "4.12", "0.0709", 5),
  @(16, "013687", "平安成长龙头1年持有混合A", "1.34", "92.94", "4.56", "0.0611", 4),
  @(17, "006100", "平安优势产业灵活配置混合A", "1.67", "93.46", "3.56", "0.0595", 8),
  @(18, "005265", "博时厚泽回报灵活配置混合A", "1.70", "72.26", "3.06", "0.0520", 7),
  @(19, "010127", "平安价值成长混合C", "1.41", "93.06", "3.63", "0.0512", 6),
  @(20, "011829", "平安睿享成长混合C", "1.14", "93.00", "4.07", "0.0464", 6),
  @(21, "000969", "前海开源大安全核心精选灵活配置混合", "1.17", "91.57", "3.96", "0.0463", 8),
  @(22, "012986", "平安优势回报1年持有期混合型证券投资基金C", "1.25", "93.07", "3.56", "0.0445", 7),
  @(23, "011885", "工银景气优选混合C", "1.03", "86.74", "2.98", "0.0307", 9),
  @(24, "005266", "博时厚泽回报灵活配置混合C", "0.96", "72.26", "3.06", "0.0294", 7),
  @(25, "009488", "中邮价值精选混合A", "0.73", "77.57", "3.90", "0.0285", 9),
  @(26, "013688", "平安成长龙头1年持有混合C", "0.56", "92.94", "4.56", "0.0255", 4),
  @(27, "002861", "工银智能制造股票", "0.68", "82.60", "3.55", "0.0241", 7),
  @(28, "007894", "平安估值精选混合C", "0.51", "92.95", "4.20", "0.0214", 6),
  @(29, "007893", "平安估值精选混合A", "0.35", "92.95", "4.20", "0.0147", 6),
  @(30, "013418", "博时核心资产精选混合C", "0.42", "70.34", "3.31", "0.0139", 6),
  @(31, "009489", "中邮价值精选混合C", "0.32", "77.57", "3.90", "0.0125", 9),
  @(32, "001721", "工银新增益混合", "0.60", "24.46", "1.00", "0.0060", 9),
  @(33, "013312", "工银主题策略混合C", "0.07", "90.81", "4.72", "0.0033", 5),
  @(34, "012918", "平安优势领航1年持有期混合C", "0.09", "93.10", "3.39", "0.0031", 8),
  @(35, "016751", "博时精选混合C", "0.00", "66.67", "2.93", "0", 7)
)

$r = 2
foreach ($row in $data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Row 37 (fund "016751") is the one exception: its "持有市值" column is a
# genuine numeric 0, not the formatted-string "0.xxxx" every other row uses.
$q3.Range("G37").NumberFormat = "General"
$q3.Cells.Item(37, 7).Value = 0

Write-Output "2022-Q3 sheet populated"
